$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "1.001") must be
# forced to Text format first, otherwise Excel auto-converts the literal
# into a numeric value instead of keeping it as the original text string.
$ws.Range("D2").Value = '23.989.74'
$ws.Range("E2").Value = '  -1.94%  '
$ws.Range("D3").Value = '1.651.08'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.10'
$ws.Range("E5").Value = '  -1.12%  '
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3908'
$ws.Range("E7").Value = '  -1.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3816'
$ws.Range("E8").Value = '  -2.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.21'
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("E10").Value = '  -4.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.001'
$ws.Range("E11").Value = '  +0.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08458'
$ws.Range("E12").Value = '  -1.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.89'
$ws.Range("E13").Value = '  -2.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.072'
$ws.Range("E14").Value = '  -3.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.007'
$ws.Range("E15").Value = '  +1.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001310'
$ws.Range("E16").Value = '  -2.79%  '
$ws.Range("D17").Value = '1.648.82'
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.30'
$ws.Range("E18").Value = '  -1.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07005'
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("E20").Value = '  -4.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.988'
$ws.Range("E21").Value = '  -0.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.80'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").Value = '23.980.77'
$ws.Range("E24").Value = '  -1.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.444'
$ws.Range("E25").Value = '  +0.61%  '
$ws.Range("E26").Value = '  -2.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.09'
$ws.Range("E27").Value = '  -1.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.99'
$ws.Range("E28").Value = '  -2.97%  '
$ws.Range("E29").Value = '  -1.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '138.14'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.934'
$ws.Range("E31").Value = '  -2.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.518'
$ws.Range("E32").Value = '  -1.12%  '
$ws.Range("D33").Value = '1.832.31'
$ws.Range("E33").Value = '  -1.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.024'
$ws.Range("E34").Value = '  -3.73%  '
$ws.Range("E35").Value = '  -2.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.770'
$ws.Range("E36").Value = '  -1.15%  '
$ws.Range("E37").Value = '  -3.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.78'
$ws.Range("E38").Value = '  -3.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2682'
$ws.Range("E39").Value = '  -3.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09106'
$ws.Range("E40").Value = '  -1.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7622'
$ws.Range("E41").Value = '  -1.93%  '
$ws.Range("E42").Value = '  -3.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.420'
$ws.Range("E43").Value = '  -1.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.30'
$ws.Range("E44").Value = '  -1.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6985'
$ws.Range("E45").Value = '  -2.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.459'
$ws.Range("E46").Value = '  -3.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.098'
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("E49").Value = '  -1.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.95'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.232'
$ws.Range("E51").Value = '  -4.09%  '
